$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1:H1").EntireColumn.Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Cells.Item(2, 6).Value = 0.0476336
$ws.Cells.Item(2, 7).Value = 0.0033069
$ws.Cells.Item(2, 8).Value = 0.0509405
$ws.Cells.Item(2, 9).Value = 9.8
$ws.Cells.Item(3, 6).Value = 0.0464272
$ws.Cells.Item(3, 7).Value = 0.00271281
$ws.Cells.Item(3, 8).Value = 0.04914
$ws.Cells.Item(3, 9).Value = 8.8
$ws.Cells.Item(4, 6).Value = 0.0459506
$ws.Cells.Item(4, 7).Value = 0.0026963
$ws.Cells.Item(4, 8).Value = 0.0486469
$ws.Cells.Item(4, 9).Value = 8.7
$ws.Cells.Item(5, 6).Value = 0.0467455
$ws.Cells.Item(5, 7).Value = 0.0032264
$ws.Cells.Item(5, 8).Value = 0.0499719
$ws.Cells.Item(5, 9).Value = 8.8
$ws.Cells.Item(6, 6).Value = 0.0468162
$ws.Cells.Item(6, 7).Value = 0.00274581
$ws.Cells.Item(6, 8).Value = 0.049562
$ws.Cells.Item(6, 9).Value = 8.5
$ws.Cells.Item(7, 6).Value = 0.0460983
$ws.Cells.Item(7, 7).Value = 0.00237766
$ws.Cells.Item(7, 8).Value = 0.0484759
$ws.Cells.Item(7, 9).Value = 8.8
$ws.Cells.Item(8, 6).Value = 0.0479079
$ws.Cells.Item(8, 7).Value = 0.00304461
$ws.Cells.Item(8, 8).Value = 0.0509525
$ws.Cells.Item(8, 9).Value = 9.8
$ws.Cells.Item(9, 6).Value = 0.048095
$ws.Cells.Item(9, 7).Value = 0.00254395
$ws.Cells.Item(9, 8).Value = 0.050639
$ws.Cells.Item(9, 9).Value = 8.3
$ws.Cells.Item(10, 6).Value = 0.0455709
$ws.Cells.Item(10, 7).Value = 0.00236888
$ws.Cells.Item(10, 8).Value = 0.0479398
$ws.Cells.Item(10, 9).Value = 8.8
$ws.Cells.Item(11, 6).Value = 0.0456859
$ws.Cells.Item(11, 7).Value = 0.0025821
$ws.Cells.Item(11, 8).Value = 0.048268
$ws.Cells.Item(11, 9).Value = 8.7
$ws.Cells.Item(12, 6).Value = 0.04669311
$ws.Cells.Item(12, 7).Value = 0.002760542
$ws.Cells.Item(12, 8).Value = 0.04945365
$ws.Cells.Item(12, 9).Value = 8.9

$ws.Range("J2:K12").ClearContents()

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)

$ws.Range("G1:H1").EntireColumn.Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Cells.Item(2, 6).Value = 0.0890898
$ws.Cells.Item(2, 7).Value = 0.00583576
$ws.Cells.Item(2, 8).Value = 0.0949255
$ws.Cells.Item(2, 9).Value = 15.7
$ws.Cells.Item(3, 6).Value = 0.0885563
$ws.Cells.Item(3, 7).Value = 0.0051259
$ws.Cells.Item(3, 8).Value = 0.0936823
$ws.Cells.Item(3, 9).Value = 15.7
$ws.Cells.Item(4, 6).Value = 0.0871721
$ws.Cells.Item(4, 7).Value = 0.00400175
$ws.Cells.Item(4, 8).Value = 0.0911739
$ws.Cells.Item(4, 9).Value = 14.5
$ws.Cells.Item(5, 6).Value = 0.0902624
$ws.Cells.Item(5, 7).Value = 0.00566894
$ws.Cells.Item(5, 8).Value = 0.0959313
$ws.Cells.Item(5, 9).Value = 15.1
$ws.Cells.Item(6, 6).Value = 0.0905172
$ws.Cells.Item(6, 7).Value = 0.00521189
$ws.Cells.Item(6, 8).Value = 0.0957291
$ws.Cells.Item(6, 9).Value = 15
$ws.Cells.Item(7, 6).Value = 0.0866788
$ws.Cells.Item(7, 7).Value = 0.00407781
$ws.Cells.Item(7, 8).Value = 0.0907566
$ws.Cells.Item(7, 9).Value = 14.5
$ws.Cells.Item(8, 6).Value = 0.0882142
$ws.Cells.Item(8, 7).Value = 0.00450065
$ws.Cells.Item(8, 8).Value = 0.0927148
$ws.Cells.Item(8, 9).Value = 15.7
$ws.Cells.Item(9, 6).Value = 0.0911633
$ws.Cells.Item(9, 7).Value = 0.00635459
$ws.Cells.Item(9, 8).Value = 0.0975178
$ws.Cells.Item(9, 9).Value = 15
$ws.Cells.Item(10, 6).Value = 0.0892155
$ws.Cells.Item(10, 7).Value = 0.0055552
$ws.Cells.Item(10, 8).Value = 0.0947707
$ws.Cells.Item(10, 9).Value = 15.1
$ws.Cells.Item(11, 6).Value = 0.0873261
$ws.Cells.Item(11, 7).Value = 0.00400778
$ws.Cells.Item(11, 8).Value = 0.0913339
$ws.Cells.Item(11, 9).Value = 15.7
$ws.Cells.Item(12, 6).Value = 0.08881957
$ws.Cells.Item(12, 7).Value = 0.005034027000000001
$ws.Cells.Item(12, 8).Value = 0.09385359
$ws.Cells.Item(12, 9).Value = 15.2

$ws.Range("J2:K12").ClearContents()

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("G1:H1").EntireColumn.Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Cells.Item(2, 6).Value = 0.132824
$ws.Cells.Item(2, 7).Value = 0.0079023
$ws.Cells.Item(2, 8).Value = 0.140727
$ws.Cells.Item(2, 9).Value = 22.1
$ws.Cells.Item(3, 6).Value = 0.134056
$ws.Cells.Item(3, 7).Value = 0.00788706
$ws.Cells.Item(3, 8).Value = 0.141943
$ws.Cells.Item(3, 9).Value = 22.1
$ws.Cells.Item(4, 6).Value = 0.138109
$ws.Cells.Item(4, 7).Value = 0.00820613
$ws.Cells.Item(4, 8).Value = 0.146315
$ws.Cells.Item(4, 9).Value = 23.2
$ws.Cells.Item(5, 6).Value = 0.132919
$ws.Cells.Item(5, 7).Value = 0.00700027
$ws.Cells.Item(5, 8).Value = 0.13992
$ws.Cells.Item(5, 9).Value = 22.4
$ws.Cells.Item(6, 6).Value = 0.136478
$ws.Cells.Item(6, 7).Value = 0.0083562
$ws.Cells.Item(6, 8).Value = 0.144834
$ws.Cells.Item(6, 9).Value = 23.6
$ws.Cells.Item(7, 6).Value = 0.144364
$ws.Cells.Item(7, 7).Value = 0.00856986
$ws.Cells.Item(7, 8).Value = 0.152934
$ws.Cells.Item(7, 9).Value = 23.6
$ws.Cells.Item(8, 6).Value = 0.139338
$ws.Cells.Item(8, 7).Value = 0.00907865
$ws.Cells.Item(8, 8).Value = 0.148417
$ws.Cells.Item(8, 9).Value = 23.3
$ws.Cells.Item(9, 6).Value = 0.137298
$ws.Cells.Item(9, 7).Value = 0.00774202
$ws.Cells.Item(9, 8).Value = 0.14504
$ws.Cells.Item(9, 9).Value = 23.6
$ws.Cells.Item(10, 6).Value = 0.138809
$ws.Cells.Item(10, 7).Value = 0.00800993
$ws.Cells.Item(10, 8).Value = 0.146819
$ws.Cells.Item(10, 9).Value = 23.3
$ws.Cells.Item(11, 6).Value = 0.138587
$ws.Cells.Item(11, 7).Value = 0.00820309
$ws.Cells.Item(11, 8).Value = 0.14679
$ws.Cells.Item(11, 9).Value = 23.3
$ws.Cells.Item(12, 6).Value = 0.1372782
$ws.Cells.Item(12, 7).Value = 0.008095551
$ws.Cells.Item(12, 8).Value = 0.1453739
$ws.Cells.Item(12, 9).Value = 23.05

$ws.Range("J2:K12").ClearContents()

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)

$ws.Range("G1:H1").EntireColumn.Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Cells.Item(2, 6).Value = 0.174549
$ws.Cells.Item(2, 7).Value = 0.0100449
$ws.Cells.Item(2, 8).Value = 0.184594
$ws.Cells.Item(2, 9).Value = 28
$ws.Cells.Item(3, 6).Value = 0.173385
$ws.Cells.Item(3, 7).Value = 0.0102257
$ws.Cells.Item(3, 8).Value = 0.18361
$ws.Cells.Item(3, 9).Value = 27.5
$ws.Cells.Item(4, 6).Value = 0.176883
$ws.Cells.Item(4, 7).Value = 0.0101787
$ws.Cells.Item(4, 8).Value = 0.187062
$ws.Cells.Item(4, 9).Value = 28.6
$ws.Cells.Item(5, 6).Value = 0.175848
$ws.Cells.Item(5, 7).Value = 0.00937562
$ws.Cells.Item(5, 8).Value = 0.185224
$ws.Cells.Item(5, 9).Value = 28.3
$ws.Cells.Item(6, 6).Value = 0.17814
$ws.Cells.Item(6, 7).Value = 0.0108777
$ws.Cells.Item(6, 8).Value = 0.189018
$ws.Cells.Item(6, 9).Value = 28
$ws.Cells.Item(7, 6).Value = 0.175022
$ws.Cells.Item(7, 7).Value = 0.0109945
$ws.Cells.Item(7, 8).Value = 0.186017
$ws.Cells.Item(7, 9).Value = 27.7
$ws.Cells.Item(8, 6).Value = 0.177476
$ws.Cells.Item(8, 7).Value = 0.0113741
$ws.Cells.Item(8, 8).Value = 0.18885
$ws.Cells.Item(8, 9).Value = 28
$ws.Cells.Item(9, 6).Value = 0.181928
$ws.Cells.Item(9, 7).Value = 0.0107158
$ws.Cells.Item(9, 8).Value = 0.192644
$ws.Cells.Item(9, 9).Value = 27.5
$ws.Cells.Item(10, 6).Value = 0.177871
$ws.Cells.Item(10, 7).Value = 0.0104599
$ws.Cells.Item(10, 8).Value = 0.188331
$ws.Cells.Item(10, 9).Value = 28
$ws.Cells.Item(11, 6).Value = 0.178217
$ws.Cells.Item(11, 7).Value = 0.0110176
$ws.Cells.Item(11, 8).Value = 0.189235
$ws.Cells.Item(11, 9).Value = 28.2
$ws.Cells.Item(12, 6).Value = 0.1769319
$ws.Cells.Item(12, 7).Value = 0.010526452
$ws.Cells.Item(12, 8).Value = 0.1874585
$ws.Cells.Item(12, 9).Value = 27.98

$ws.Range("J2:K12").ClearContents()
